# Generate Report for Handback
# ---------------------------------------------------------------
# This applies the "handback" localization-status report update:
#   * Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#   * Stamps a "Latest Handback DateTime" for each locale row
#   * Fills in the (previously empty) "Latest Target File" / "Latest
#     Handback File" columns, turning "Latest Target File" into a
#     hyperlink (same md file as the Source File Name column)
#   * Widens a few columns that now hold longer content
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$mdUrlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dae4197392d746f947c8871cb467bf5b6d498dc3/e2e/"
$file1 = "4671043e-79a3-44df-ba54-b798b1604ef4"
$file2 = "5a546f92-ab30-4fcf-bbec-bfe5f7c73a43"
$file1Md = "$file1.md"
$file2Md = "$file2.md"

function Update-LocaleSheet($SheetName, $Locale, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # --- Status column (C) -> handed back ---
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # --- Latest Handback DateTime (K) ---
    $ws.Range("K2").Value = $HandbackDateTime
    $ws.Range("K3").Value = $HandbackDateTime

    # --- Latest Handback File (J) - plain text, matches the xlf already
    #     referenced under "Latest Handoff File" (G) for this locale ---
    $ws.Range("J2").Value = "$file1.4ddb0be76cb358a7c0ce4470de5500fd3755ce1d.$Locale.xlf"
    $ws.Range("J3").Value = "$file2.d9abb552313603295ad77449c1f87efd3a2edee5.$Locale.xlf"

    # --- Latest Target File (I) - hyperlink to the same .md source doc
    #     used by column A, rebuild all 4 hyperlinks in display order so
    #     relationship ids line up (A2, I2, A3, I3) ---
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "$mdUrlBase$file1Md", "", "", $file1Md)
    $ws.Hyperlinks.Add($ws.Range("I2"), "$mdUrlBase$file1Md", "", "", $file1Md)
    $ws.Hyperlinks.Add($ws.Range("A3"), "$mdUrlBase$file2Md", "", "", $file2Md)
    $ws.Hyperlinks.Add($ws.Range("I3"), "$mdUrlBase$file2Md", "", "", $file2Md)

    # --- column widths for the now-wider columns ---
    $ws.Columns.Item(3).ColumnWidth = 29.1
    $ws.Columns.Item(9).ColumnWidth = 39.1
    $ws.Columns.Item(10).ColumnWidth = 39.1
}

Update-LocaleSheet "zh-cn" "zh-cn" "2016-09-04 18:52:55"
Update-LocaleSheet "de-de" "de-de" "2016-09-04 18:53:07"

# --- Overview sheet: Status columns (E, F) just widen along with the
#     text change that already propagated via the shared string update
#     above ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

Write-Host "Handback report generated."
